# Fix string/int error in the Year_Range labels and update the shuffled
# Keyword/Correlation values for the 2020-01_2022-05 and 2022-06_2023-12
# periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Year_Range text in column A -------------------------------
$ws.Range("A2:A6").Value = "2007-12_2009-06"
$ws.Range("A7:A11").Value = "2009-07_2019-12"
$ws.Range("A12:A16").Value = "2020-01_2022-05"
$ws.Range("A17:A21").Value = "2022-06_2023-12"

# --- Rows 12-16 (2020-01_2022-05 block): Keyword + Correlation reshuffle
$ws.Range("B12").Value = "inflation"
$ws.Range("C12").Value = 0.821

$ws.Range("B13").Value = "uncertain"
$ws.Range("C13").Value = 0.0716

$ws.Range("B14").Value = "interest"
$ws.Range("C14").Value = 0.8708

$ws.Range("B15").Value = "invest"
$ws.Range("C15").Value = 0.1018

$ws.Range("B16").Value = "trade"
$ws.Range("C16").Value = 0.706

# --- Rows 17-21 (2022-06_2023-12 block): Correlation updates only ------
$ws.Range("C17").Value = -0.8116
$ws.Range("C18").Value = -0.2116
$ws.Range("C19").Value = 0.2492
$ws.Range("C20").Value = 0.6523
$ws.Range("C21").Value = 0.4977
